$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C4/D4: they were stored as text, now should be real numbers
$ws.Range("C4").Value = 123456789
$ws.Range("D4").Value = 12345678

# Add a new row of data for client "amo" / "Amos"
$ws.Range("A5").Value = "amo"
$ws.Range("B5").Value = "Amos"

# Contact / IFU are left blank for this client, but still present as
# (text-typed) empty cells on the row, matching the other data rows.
$ws.Range("C5").Value = "'"
$ws.Range("D5").Value = "'"
